$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source")

# Update data values on the "Source" sheet
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("I2").Value = 8

$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 15
$ws.Range("G3").Value = 15
$ws.Range("I3").Value = 15

$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 0.5
$ws.Range("G4").Value = 0.5
$ws.Range("I4").Value = 0.5

# Update selected cell / window view
$ws.Range("E14").Select()

# Reflect the enlarged application window size recorded in the saved file
$excel.ActiveWindow.Width = 20490
$excel.ActiveWindow.Height = 7545
